# Update the "想去人数" (F column) figures in the 展览 and 全部类型 sheets.
# These two sheets carry duplicate data, so the same row/value updates
# are applied to both.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 2730
    4  = 613
    6  = 6639
    7  = 1119
    8  = 12
    9  = 18
    10 = 16
    11 = 67
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
